$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column R: year 2021 ---
# Header cell R4: copy formatting from the existing 2020 header (Q4), then set the new year value.
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

# Data cell R5: copy formatting from the existing 2020 data cell (Q5), then tweak number
# format / alignment for the new value and set it.
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").NumberFormat = "0.0"
$ws.Range("R5").HorizontalAlignment = -4152
$ws.Range("R5").Value = 102.20441221981518

# Restore the selection to match the state after the edit.
[void]$ws.Range("S9").Select()
